# Generate Report for Handback
#
# The 2f0ceece-... file's handback failed because the returned file name
# didn't match the original handoff file name. Reflect that in the
# status columns and record the error detail for both locales.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update status from "Ready for handoff" to "Handback transform failed"
# for the 2f0ceece-... row, everywhere it is shown.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback error detail for each locale and widen the
# "Error Detail" column (P) so the message is readable.
$wsZhCn.Range("P3").Value = "Handback file name: tgpjxfr4.hh1 is different with handoff file name: 2f0ceece-0e35-42cd-8e92-33eb2ea97769.9dbe452733378811f4c516f73ad8ac6db13439a4.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: tgpjxfr4.hh1 is different with handoff file name: 2f0ceece-0e35-42cd-8e92-33eb2ea97769.9dbe452733378811f4c516f73ad8ac6db13439a4.de-de."

# Excel pads ColumnWidth by ~5/6 of a character (the "Normal style" font
# metric) when it stores the width in the OOXML <col> element, so dial the
# requested width back by that amount to land on a stored width of 40.
$colWidth = 40 - (5 / 6)
$wsZhCn.Columns.Item(16).ColumnWidth = $colWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $colWidth
